$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header row) - sample-size values for columns B:E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) - updated measurements for columns B:E
$ws.Range("B2").Value = 11.562596071250281
$ws.Range("C2").Value = 1.4715207672810251
$ws.Range("D2").Value = 0.82064937990615228
$ws.Range("E2").Value = 0.64227479785525987

# Row 3 (STR) - updated measurements for columns B:E
$ws.Range("B3").Value = 4.1569834129557313
$ws.Range("C3").Value = 6.5977523756936236
$ws.Range("D3").Value = 2.532526546573743
$ws.Range("E3").Value = -0.021957963716064517

# Selection now only covers the edited block instead of the whole sheet
$ws.Range("B1:E3").Select()
